# Applies the changes described by the diff to the Functional Requirement
# Document:
#   1. The empty paragraph after the "Listings can have different
#      visibility settings..." bullet is turned into two new bulleted
#      paragraphs about verified community members posting services.
#   2/3/4/5. Four "lastRenderedPageBreak" markers shift position because of
#      the reflow caused by the new content (removed from one run, added to
#      another, for three separate spots in the document).

$d = $word.ActiveDocument
$wdParagraph = 4
$wdCollapseEnd = 0
$wmain = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Get-ParagraphRangeAfter($searchText) {
    # Returns a Range covering the paragraph immediately following the one
    # containing $searchText.
    $rng = $d.Content
    $null = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $rng.Expand($wdParagraph) | Out-Null
    $rng.Collapse($wdCollapseEnd) | Out-Null
    $rng.Expand($wdParagraph) | Out-Null
    return $rng
}

function Get-ParagraphRange($searchText) {
    # Returns a Range covering the whole paragraph containing $searchText.
    $rng = $d.Content
    $null = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $rng.Expand($wdParagraph) | Out-Null
    return $rng
}

# --- 1. Empty paragraph -> two new bulleted list items ---------------------
$target = Get-ParagraphRangeAfter("Listings can have different visibility settings")
$xml = '<w:p xmlns:w="' + $wmain + '">' +
         '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
         '<w:r><w:lastRenderedPageBreak/><w:t>Verified community member</w:t></w:r>' +
         '<w:r><w:t>s</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> can post the servic</w:t></w:r>' +
         '<w:r><w:t>e by buying the service posting license</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> through platform.</w:t></w:r>' +
       '</w:p>' +
       '<w:p xmlns:w="' + $wmain + '">' +
         '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
         '<w:r><w:t>Anyone who is not a verified community member needs</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> to contact us, buy the license and then can post the service.</w:t></w:r>' +
       '</w:p>'
$target.InsertXML($xml)

# --- 2. Remove the lastRenderedPageBreak before "3. User Registration..." --
$target = Get-ParagraphRange("3. User Registration")
$xml = '<w:p xmlns:w="' + $wmain + '" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4B05990A" w14:textId="3AD7BEF8" w:rsidR="00CB6ADD" w:rsidRPr="00046320" w:rsidRDefault="00CB6ADD" w:rsidP="00CB6ADD">' +
         '<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>3</w:t></w:r>' +
         '<w:r w:rsidRPr="00046320"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>. User Registration &amp; Authentication</w:t></w:r>' +
       '</w:p>'
$target.InsertXML($xml)

# --- 3. Add lastRenderedPageBreak before "He cannot share his interest..." -
$target = Get-ParagraphRange("He cannot share his interest with the seller.")
$xml = '<w:p xmlns:w="' + $wmain + '" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7D9FA0C6" w14:textId="77777777" w:rsidR="003A36C5" w:rsidRDefault="007312B3" w:rsidP="0077105F">' +
         '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="16"/></w:numPr></w:pPr>' +
         '<w:r w:rsidRPr="007312B3"><w:lastRenderedPageBreak/><w:t>He cannot share his interest with the seller.</w:t></w:r>' +
       '</w:p>'
$target.InsertXML($xml)

# --- 4. Remove lastRenderedPageBreak before "He can also verify..." --------
$target = Get-ParagraphRange("He can also verify him/herself by")
$xml = '<w:p xmlns:w="' + $wmain + '" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="16D0DD0A" w14:textId="07967029" w:rsidR="007312B3" w:rsidRPr="007312B3" w:rsidRDefault="007312B3" w:rsidP="0064086B">' +
         '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr>' +
         '<w:r w:rsidRPr="007312B3"><w:t xml:space="preserve">He can also verify him/herself by </w:t></w:r>' +
         '<w:r w:rsidR="00D1035B" w:rsidRPr="007312B3"><w:t>uploading</w:t></w:r>' +
         '<w:r w:rsidRPr="007312B3"><w:t xml:space="preserve"> the resident proof and become a Verified community member. </w:t></w:r>' +
       '</w:p>'
$target.InsertXML($xml)

# --- 5. Add lastRenderedPageBreak before "6. He can view the conversations" -
$target = Get-ParagraphRange("6. He can view the conversations")
$xml = '<w:p xmlns:w="' + $wmain + '" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="24E86575" w14:textId="77777777" w:rsidR="007312B3" w:rsidRPr="007312B3" w:rsidRDefault="007312B3" w:rsidP="007312B3">' +
         '<w:r w:rsidRPr="007312B3"><w:lastRenderedPageBreak/><w:t>        6. He can view the conversations with other verified community members or</w:t></w:r>' +
       '</w:p>'
$target.InsertXML($xml)

# --- 6. Remove lastRenderedPageBreak before "7. He can post an item..." ----
$target = Get-ParagraphRange("7. He can post an item for sale")
$xml = '<w:p xmlns:w="' + $wmain + '" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0DD2B939" w14:textId="77777777" w:rsidR="005F1498" w:rsidRDefault="007312B3" w:rsidP="007312B3">' +
         '<w:r w:rsidRPr="007312B3"><w:t>        7. He can post an item for sale(although not a new item)</w:t></w:r>' +
       '</w:p>'
$target.InsertXML($xml)

Write-Output "done"
